$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 180.42857
$ws.Range("I39").Value = 43.833332
$ws.Range("K39").Value = 131.499996
$ws.Range("M39").Value = 164.500004

$ws.Range("H40").Value = 2458.3333
$ws.Range("I40").Value = 2187.5
$ws.Range("K40").Value = 2187.5
$ws.Range("M40").Value = -2012.5

$ws.Range("H98").Value = 1882.9474
$ws.Range("I98").Value = 2040.0588
$ws.Range("K98").Value = 2040.0588
$ws.Range("M98").Value = -542.0588

$ws.Range("H121").Value = 13649
$ws.Range("J121").Value = 13649
$ws.Range("L121").Value = 40947
$ws.Range("N121").Value = -44441

$ws.Range("H122").Value = 1882.9474
$ws.Range("I122").Value = 2040.0588
$ws.Range("K122").Value = 6120.1764
$ws.Range("M122").Value = -3670.1764

$ws.Range("H131").Value = 2224.125
$ws.Range("I131").Value = 1960.4615
$ws.Range("K131").Value = 5881.3845
$ws.Range("M131").Value = -841.3845000000001

$ws.Range("H135").Value = 785.2
$ws.Range("I135").Value = 766.08
$ws.Range("K135").Value = 6894.72
$ws.Range("M135").Value = -4359.72

$ws.Range("H137").Value = 1381.0435
$ws.Range("I137").Value = 1251.3334
$ws.Range("K137").Value = 3754.0002
$ws.Range("M137").Value = -1204.0002

$ws.Range("H138").Value = 4967.1606
$ws.Range("J138").Value = 7586.6216
$ws.Range("L138").Value = 22759.8648
$ws.Range("N138").Value = -33039.8648

$ws.Range("H141").Value = 4152.16
$ws.Range("I141").Value = 2126.5264
$ws.Range("K141").Value = 6379.5792
$ws.Range("M141").Value = -1199.5792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 81.5
$ws.Range("I5").Value = 20
$ws.Range("K5").Value = 20
$ws.Range("M5").Value = 92

$ws.Range("H32").Value = 5776.6445
$ws.Range("I32").Value = 3027.7896
$ws.Range("K32").Value = 3027.7896
$ws.Range("M32").Value = -2740.7896

$ws.Range("H45").Value = 1832.75
$ws.Range("I45").Value = 1817.6364
$ws.Range("K45").Value = 1817.6364
$ws.Range("M45").Value = -1440.6364

$ws.Range("H55").Value = 35713.285
$ws.Range("J55").Value = 38998.6
$ws.Range("L55").Value = 38998.6
$ws.Range("N55").Value = -39628.6

$ws.Range("H122").Value = 1712
$ws.Range("I122").Value = 1712
$ws.Range("K122").Value = 5136
$ws.Range("M122").Value = -2686

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 81.5
$ws.Range("I4").Value = 20
$ws.Range("K4").Value = 20
$ws.Range("M4").Value = 95

$ws.Range("H26").Value = 15349.75
$ws.Range("I26").Value = 15349.75
$ws.Range("K26").Value = 15349.75
$ws.Range("M26").Value = -15057.75

$ws.Range("H30").Value = 3000
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4391.3335
$ws.Range("I16").Value = 4309.3076
$ws.Range("J16").Value = 4604.6
$ws.Range("K16").Value = 4309.3076
$ws.Range("L16").Value = 4604.6
$ws.Range("M16").Value = -4022.3076
$ws.Range("N16").Value = -5178.6

$ws.Range("H31").Value = 4820.56
$ws.Range("I31").Value = 4615.4
$ws.Range("J31").Value = 5128.3
$ws.Range("K31").Value = 4615.4
$ws.Range("L31").Value = 5128.3
$ws.Range("M31").Value = -4320.4
$ws.Range("N31").Value = -5718.3

$ws.Range("H34").Value = 4820.56
$ws.Range("I34").Value = 4615.4
$ws.Range("J34").Value = 5128.3
$ws.Range("K34").Value = 4615.4
$ws.Range("L34").Value = 5128.3
$ws.Range("M34").Value = -4413.4
$ws.Range("N34").Value = -5532.3

$ws.Range("H60").Value = 15193.75
$ws.Range("J60").Value = 49997
$ws.Range("L60").Value = 49997
$ws.Range("N60").Value = -51019

$ws.Range("H113").Value = 4391.3335
$ws.Range("I113").Value = 4309.3076
$ws.Range("J113").Value = 4604.6
$ws.Range("K113").Value = 4309.3076
$ws.Range("L113").Value = 4604.6
$ws.Range("M113").Value = -2139.3076
$ws.Range("N113").Value = -8944.6

$ws.Range("H132").Value = 756.5217
$ws.Range("I132").Value = 742.2857
$ws.Range("K132").Value = 2226.8571
$ws.Range("M132").Value = 303.1428999999998

$ws.Range("H134").Value = 2138.2354
$ws.Range("I134").Value = 1148.1111
$ws.Range("K134").Value = 3444.3333
$ws.Range("M134").Value = -909.3333000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 11999.857
$ws.Range("J92").Value = 10399.8
$ws.Range("L92").Value = 10399.8
$ws.Range("N92").Value = -14143.8

$ws.Range("H122").Value = 64721.125
$ws.Range("J122").Value = 168661
$ws.Range("L122").Value = 505983
$ws.Range("N122").Value = -510883

$ws.Range("H123").Value = 24461.77
$ws.Range("J123").Value = 24461.77
$ws.Range("L123").Value = 24461.77
$ws.Range("N123").Value = -29361.77

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3652
$ws.Range("I7").Value = 3478
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 3478
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -3366
$ws.Range("N7").Value = -4224

$ws.Range("H41").Value = 20000
$ws.Range("I41").Value = 20000
$ws.Range("K41").Value = 20000
$ws.Range("M41").Value = -19562

$ws.Range("H46").Value = 3269.111
$ws.Range("I46").Value = 2633.3333
$ws.Range("J46").Value = 3587
$ws.Range("K46").Value = 2633.3333
$ws.Range("L46").Value = 3587
$ws.Range("M46").Value = -2445.3333
$ws.Range("N46").Value = -3963

$ws.Range("H126").Value = 3652
$ws.Range("I126").Value = 3478
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 10434
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -7964
$ws.Range("N126").Value = -16940

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 4193.3
$ws.Range("I132").Value = 3617.25
$ws.Range("K132").Value = 10851.75
$ws.Range("M132").Value = -8321.75

$ws.Range("H136").Value = 3118
$ws.Range("I136").Value = 2699.889
$ws.Range("K136").Value = 8099.667
$ws.Range("M136").Value = -5549.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 14993.875
$ws.Range("J18").Value = 14996.8
$ws.Range("L18").Value = 14996.8
$ws.Range("N18").Value = -15342.8

$ws.Range("H70").Value = 90000
$ws.Range("I70").Value = 90000
$ws.Range("K70").Value = 90000
$ws.Range("M70").Value = -89685

$ws.Range("H73").Value = 90000
$ws.Range("I73").Value = 90000
$ws.Range("K73").Value = 90000
$ws.Range("M73").Value = -88908

$ws.Range("H93").Value = 129999
$ws.Range("J93").Value = 129999
$ws.Range("L93").Value = 129999
$ws.Range("N93").Value = -134991

$ws.Range("H107").Value = 1220.7858
$ws.Range("I107").Value = 1254.2222
$ws.Range("J107").Value = 1160.6
$ws.Range("K107").Value = 3762.6666
$ws.Range("L107").Value = 3481.8
$ws.Range("M107").Value = -1842.6666
$ws.Range("N107").Value = -7321.799999999999

$ws.Range("H126").Value = 2042.3334
$ws.Range("I126").Value = 1368.5
$ws.Range("J126").Value = 3390
$ws.Range("K126").Value = 4105.5
$ws.Range("L126").Value = 10170
$ws.Range("M126").Value = -1635.5
$ws.Range("N126").Value = -15110

$ws.Range("H132").Value = 41982.566
$ws.Range("I132").Value = 55918.117
$ws.Range("K132").Value = 167754.351
$ws.Range("M132").Value = -165224.351

$ws.Range("H135").Value = 69999
$ws.Range("J135").Value = 69999
$ws.Range("L135").Value = 69999
$ws.Range("N135").Value = -80139
